# Rotate values in columns B, D, L, O, P across rows 8, 9, 10:
#   new row8  <- old row10
#   new row9  <- old row8
#   new row10 <- old row9
#
# Notes on this host's COM surface:
#  - Reading the plain `.Value` getter on a Range does not resolve to the
#    cell's actual contents here -- use `.Value2` (or `.Text`/`.Formula`)
#    for reads. `.Value` as a setter works fine, but for symmetry this
#    script reads and writes via `.Value2` throughout.
#  - Writing an empty string ("") always blanks the cell out completely
#    (same as genuine Excel: there is no Range.Value assignment that
#    leaves a cell holding a zero-length string instead of being empty).
#    So cells that are not actually changing value are skipped, to avoid
#    turning an already-blank cell into a "different kind of blank".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "L", "O", "P")

foreach ($col in $cols) {
    $v8  = $ws.Range("$col" + "8").Value2
    $v9  = $ws.Range("$col" + "9").Value2
    $v10 = $ws.Range("$col" + "10").Value2

    if ($v8 -ne $v10) { $ws.Range("$col" + "8").Value2  = $v10 }
    if ($v9 -ne $v8)  { $ws.Range("$col" + "9").Value2  = $v8 }
    if ($v10 -ne $v9) { $ws.Range("$col" + "10").Value2 = $v9 }
}
